# Weekly update: a new price observation was inserted as row 242
# ("Haba", origin Carahue, fecha 44798) and every subsequent record on
# this sheet shifted down by one row (old row 242 -> new row 243, ...,
# old row 294 -> new row 295).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 242:294 down to 243:295, opening up a blank row 242.
$ws.Rows(242).Insert()

# Populate the newly inserted row with the new observation.
$ws.Range("A242").Value = 6
$ws.Range("B242").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C242").Value = "Metropolitana"
$ws.Range("D242").Value = 44798
$ws.Range("E242").Value = 13
$ws.Range("F242").Value = 100112026
$ws.Range("G242").Value = "Haba"
$ws.Range("H242").Value = "Sin especificar"
$ws.Range("I242").Value = "Primera"
$ws.Range("J242").Value = 400
$ws.Range("K242").Value = 17000
$ws.Range("L242").Value = 19000
$ws.Range("M242").Value = 17850
$ws.Range("N242").Value = "$/saco 25 kilos"
$ws.Range("O242").Value = "Carahue"
$ws.Range("P242").Value = 714
$ws.Range("Q242").Value = 25
$ws.Range("R242").Value = "Hortaliza"
